# CRUD para alimentação do programa implementado
# Update the "agendamentos" sheet:
#  - Row 2 (abraao / Paulo / Pintura @ 26/05/2026 20:00)
#  - Row 3 (abraao / João  / Pintura @ 25/01/2026 15:00)
#  - Remove the former row 4 (abraao / Pacote Completo) entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 4) first, so indices of remaining rows stay 2 and 3
$ws.Rows.Item(4).Delete()

# Row 2
$ws.Range("A2").Value = "abraao"
$ws.Range("B2").Value = "abraaocursos2019@gmail.com"
$ws.Range("C2").Value = 85986820652
$ws.Range("D2").Value = "Pintura"
$ws.Range("E2").Value = "Pintura de cabelo com tinta temporária"
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = "Paulo"
$ws.Range("H2").Value = "26/05/2026"
$ws.Range("I2").Value = "20:00"
$ws.Range("J2").Value = "2026-05-26T20:00:00"
$ws.Range("K2").Value = "agendamento"

# Row 3
$ws.Range("A3").Value = "abraao"
$ws.Range("B3").Value = "abraaocursos2019@gmail.com"
$ws.Range("C3").Value = 85986820652
$ws.Range("D3").Value = "Pintura"
$ws.Range("E3").Value = "Pintura de cabelo com tinta temporária"
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = "João"
$ws.Range("H3").Value = "25/01/2026"
$ws.Range("I3").Value = "15:00"
$ws.Range("J3").Value = "2026-01-25T15:00:00"
$ws.Range("K3").Value = "agendamento"
